$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new example sentence + blocks, and a plain numeric "correct index" in F2
$ws.Range("A2").Value = "This is a  ___ sentence"
$ws.Range("B2").Value = "false1"
$ws.Range("C2").Value = "false2"
$ws.Range("D2").Value = "correct"
$ws.Range("E2").Value = "false3"

# F2 used to be a styled number (s="2"); it becomes a plain, unstyled number
$ws.Range("F2").ClearFormats()
$ws.Range("F2").Value = 3

# Row 3: blank out A:E (keep the existing style), drop F3 entirely (content + format)
$ws.Range("A3:E3").ClearContents()
$ws.Range("F3").Clear()

# Rows 4 and 5: blank out all six columns but keep their existing style
$ws.Range("A4:F4").ClearContents()
$ws.Range("A5:F5").ClearContents()

# Selection moves from H6 to F2
$ws.Range("F2").Select()
